$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine how many rows of data currently exist (header + data rows)
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column before column A; this shifts the existing
# item_name/calories columns from A/B to B/C, carrying their widths,
# styles and data along with them.
$ws.Range("A1").EntireColumn.Insert()

# New column A header
$ws.Range("A1").Value = "restaurant"

# Every data row belongs to the Wendy's menu
$ws.Range("A2:A" + $lastRow).Value = "Wendys"

# Restore the (arbitrary) active cell/selection recorded in the saved view
$ws.Range("I12").Select()
